$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.425.07"
$ws.Range("E2").Value = "  +0.58%  "

$ws.Range("D3").Value = "2.246.92"
$ws.Range("E3").Value = "  -0.25%  "

$ws.Range("E4").Value = "  +0.29%  "

$ws.Range("D5").Value = "307.51"
$ws.Range("E5").Value = "  -0.21%  "

$ws.Range("D6").Value = "94.37"
$ws.Range("E6").Value = "  -4.29%  "

$ws.Range("E7").Value = "  -0.89%  "

$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.22%  "

$ws.Range("E9").Value = "  -1.44%  "

$ws.Range("D10").Value = "34.79"
$ws.Range("E10").Value = "  -2.21%  "

$ws.Range("D11").Value = "0.0812"
$ws.Range("E11").Value = "  -1.20%  "

$ws.Range("D12").Value = "7.19"
$ws.Range("E12").Value = "  -1.70%  "

$ws.Range("E13").Value = "  +0.29%  "

$ws.Range("D14").Value = "2.355.97"
$ws.Range("E14").Value = "  +4.70%  "

$ws.Range("D15").Value = "0.840"
$ws.Range("E15").Value = "  +0.15%  "

$ws.Range("D16").Value = "13.66"
$ws.Range("E16").Value = "  -1.10%  "

$ws.Range("D17").Value = "44.056.26"
$ws.Range("E17").Value = "  +0.14%  "

$ws.Range("D18").Value = "12.52"
$ws.Range("E18").Value = "  -3.75%  "

$ws.Range("D19").Value = "0.0₃0963"
$ws.Range("E19").Value = "  -1.13%  "

$ws.Range("D20").Value = "6.41"
$ws.Range("E20").Value = "  +1.41%  "

$ws.Range("D21").Value = "65.88"
$ws.Range("E21").Value = "  +0.84%  "

$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "237.50"
$ws.Range("E22").Value = "  -2.14%  "

$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").Value = "3.00"
$ws.Range("E23").Value = "  +1.70%  "

$ws.Range("D24").Value = "2.00"
$ws.Range("E24").Value = "  +1.50%  "

$ws.Range("E25").Value = "  -0.20%  "

$ws.Range("D26").Value = "38.59"
$ws.Range("E26").Value = "  +5.49%  "

$ws.Range("E27").Value = "  +3.47%  "

$ws.Range("D28").Value = "9.86"
$ws.Range("E28").Value = "  -2.42%  "

$ws.Range("E29").Value = "  -4.07%  "

$ws.Range("D30").Value = "20.11"
$ws.Range("E30").Value = "  -0.28%  "

$ws.Range("D31").Value = "154.29"
$ws.Range("E31").Value = "  -1.69%  "

$ws.Range("D32").Value = "0.0800"
$ws.Range("E32").Value = "  -3.01%  "

$ws.Range("E33").Value = "  -0.67%  "

$ws.Range("D34").Value = "3.11"
$ws.Range("E34").Value = "  -11.11%  "

$ws.Range("E35").Value = "  +2.08%  "

$ws.Range("E36").Value = "  +0.35%  "

$ws.Range("D37").Value = "1.82"
$ws.Range("E37").Value = "  -1.41%  "

$ws.Range("D38").Value = "3.46"
$ws.Range("E38").Value = "  +2.83%  "

$ws.Range("E39").Value = "  -5.55%  "

$ws.Range("D40").Value = "3.81"
$ws.Range("E40").Value = "  -1.48%  "

$ws.Range("D41").Value = "0.0305"
$ws.Range("E41").Value = "  -0.82%  "

$ws.Range("E42").Value = "  +0.35%  "

$ws.Range("D43").Value = "1.743.96"
$ws.Range("E43").Value = "  -1.41%  "

$ws.Range("E44").Value = "  +0.12%  "

$ws.Range("D45").Value = "80.80"
$ws.Range("E45").Value = "  -8.12%  "

$ws.Range("D46").Value = "70.98"
$ws.Range("E46").Value = "  +1.73%  "

$ws.Range("D47").Value = "99.62"
$ws.Range("E47").Value = "  -1.64%  "

$ws.Range("E48").Value = "  -4.13%  "

$ws.Range("D49").Value = "56.09"
$ws.Range("E49").Value = "  +0.24%  "

$ws.Range("D50").Value = "8.18"
$ws.Range("E50").Value = "  -0.67%  "

$ws.Range("E51").Value = "  +3.20%  "
